# Weekly fruit/vegetable price update: insert two new daily-price rows
# (date 2022-04-13 / serial 44664) at the top of the data block, pushing
# the rest of the historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 130 (existing rows 130.. shift down to 132..)
$ws.Rows(130).Insert()
$ws.Rows(130).Insert()

# New row 130: Primera quality
$ws.Range("A130").Value = 7
$ws.Range("B130").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C130").Value = "Ñuble"
$ws.Range("D130").Value = 44664
$ws.Range("E130").Value = 16
$ws.Range("F130").Value = 100112008
$ws.Range("G130").Value = "Coliflor"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 200
$ws.Range("K130").Value = 1100
$ws.Range("L130").Value = 1200
$ws.Range("M130").Value = 1150
$ws.Range("N130").Value = "`$/unidad"
$ws.Range("O130").Value = "Región del Maule"
$ws.Range("P130").Value = 1150
$ws.Range("Q130").Value = 1
$ws.Range("R130").Value = "Hortaliza"

# New row 131: Segunda quality
$ws.Range("A131").Value = 7
$ws.Range("B131").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C131").Value = "Ñuble"
$ws.Range("D131").Value = 44664
$ws.Range("E131").Value = 16
$ws.Range("F131").Value = 100112008
$ws.Range("G131").Value = "Coliflor"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Segunda"
$ws.Range("J131").Value = 100
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 900
$ws.Range("M131").Value = 900
$ws.Range("N131").Value = "`$/unidad"
$ws.Range("O131").Value = "Región del Maule"
$ws.Range("P131").Value = 900
$ws.Range("Q131").Value = 1
$ws.Range("R131").Value = "Hortaliza"
